$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the diff: Price (D) and Volume(1h) (E) columns refreshed.
# Some Price values are plain-decimal strings (e.g. "1.001", "1.180") that Excel
# would otherwise auto-convert to numbers (losing trailing zeros / switching to
# scientific notation). For those we force a Text quote-prefix, write the value,
# then restore the cell style so no stray number-format styling is left behind.

$ws.Range("D2").Value = '29.494.52'
$ws.Range("E2").Value = '  +3.89%  '
$ws.Range("D3").Value = '1.911.55'
$ws.Range("E3").Value = '  +2.45%  '
$c = $ws.Range("D4")
$s = $c.Style
$c.Value = "'1.001"
$c.Style = $s
$ws.Range("E4").Value = '  -0.09%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'332.71"
$c.Style = $s
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("E6").Value = '  -0.03%  '
$c = $ws.Range("D7")
$s = $c.Style
$c.Value = "'0.4673"
$c.Style = $s
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("E8").Value = '  +2.46%  '
$c = $ws.Range("D9")
$s = $c.Style
$c.Value = "'47.88"
$c.Style = $s
$ws.Range("E9").Value = '  +0.19%  '
$c = $ws.Range("D10")
$s = $c.Style
$c.Value = "'0.08025"
$c.Style = $s
$ws.Range("E10").Value = '  +2.18%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'1.012"
$c.Style = $s
$ws.Range("E11").Value = '  +2.71%  '
$c = $ws.Range("D12")
$s = $c.Style
$c.Value = "'22.36"
$c.Style = $s
$ws.Range("E12").Value = '  +5.04%  '
$ws.Range("D13").Value = '1.909.47'
$ws.Range("E13").Value = '  +2.18%  '
$c = $ws.Range("D14")
$s = $c.Style
$c.Value = "'5.975"
$c.Style = $s
$c = $ws.Range("D15")
$s = $c.Style
$c.Value = "'7.175"
$c.Style = $s
$ws.Range("E15").Value = '  +2.66%  '
$c = $ws.Range("D16")
$s = $c.Style
$c.Value = "'89.83"
$c.Style = $s
$ws.Range("E16").Value = '  +1.88%  '
$c = $ws.Range("D17")
$s = $c.Style
$c.Value = "'1.001"
$c.Style = $s
$ws.Range("E17").Value = '  -0.06%  '
$c = $ws.Range("D18")
$s = $c.Style
$c.Value = "'0.00001032"
$c.Style = $s
$ws.Range("E18").Value = '  +1.34%  '
$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'0.06592"
$c.Style = $s
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '29.463.27'
$ws.Range("E22").Value = '  +3.83%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.Value = "'5.571"
$c.Style = $s
$ws.Range("E23").Value = '  +4.19%  '
$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'11.51"
$c.Style = $s
$ws.Range("E24").Value = '  +5.97%  '
$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'2.214"
$c.Style = $s
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").Value = '2.145.51'
$ws.Range("E26").Value = '  +2.63%  '
$c = $ws.Range("D27")
$s = $c.Style
$c.Value = "'154.74"
$c.Style = $s
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("E28").Value = '  +2.70%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.Value = "'5.759"
$c.Style = $s
$ws.Range("E29").Value = '  +8.78%  '
$c = $ws.Range("D30")
$s = $c.Style
$c.Value = "'2.141"
$c.Style = $s
$ws.Range("E30").Value = '  +4.13%  '
$c = $ws.Range("D31")
$s = $c.Style
$c.Value = "'117.30"
$c.Style = $s
$ws.Range("E31").Value = '  -0.20%  '
$c = $ws.Range("D32")
$s = $c.Style
$c.Value = "'1.063"
$c.Style = $s
$ws.Range("E32").Value = '  +10.83%  '
$c = $ws.Range("D33")
$s = $c.Style
$c.Value = "'0.09461"
$c.Style = $s
$ws.Range("E33").Value = '  +1.23%  '
$c = $ws.Range("D34")
$s = $c.Style
$c.Value = "'1.426"
$c.Style = $s
$ws.Range("E34").Value = '  +2.77%  '
$c = $ws.Range("D35")
$s = $c.Style
$c.Value = "'3.573"
$c.Style = $s
$ws.Range("E35").Value = '  -0.17%  '
$c = $ws.Range("D36")
$s = $c.Style
$c.Value = "'5.403"
$c.Style = $s
$ws.Range("E36").Value = '  +2.98%  '
$c = $ws.Range("D37")
$s = $c.Style
$c.Value = "'0.06115"
$c.Style = $s
$ws.Range("E37").Value = '  +1.32%  '
$c = $ws.Range("D38")
$s = $c.Style
$c.Value = "'0.02262"
$c.Style = $s
$ws.Range("E38").Value = '  +2.74%  '
$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'8.406"
$c.Style = $s
$ws.Range("E39").Value = '  +1.35%  '
$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'1.180"
$c.Style = $s
$ws.Range("E40").Value = '  +1.68%  '
$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'0.5888"
$c.Style = $s
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("E43").Value = '  +1.43%  '
$c = $ws.Range("D44")
$s = $c.Style
$c.Value = "'1.270"
$c.Style = $s
$ws.Range("E44").Value = '  -0.01%  '
$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'2.352"
$c.Style = $s
$ws.Range("E45").Value = '  +2.90%  '
$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'0.07509"
$c.Style = $s
$ws.Range("E46").Value = '  +5.50%  '
$c = $ws.Range("D47")
$s = $c.Style
$c.Value = "'0.5564"
$c.Style = $s
$ws.Range("E47").Value = '  +2.38%  '
$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'12.18"
$c.Style = $s
$ws.Range("E48").Value = '  +1.99%  '
$c = $ws.Range("D49")
$s = $c.Style
$c.Value = "'1.926"
$c.Style = $s
$ws.Range("E49").Value = '  +2.07%  '
$c = $ws.Range("D50")
$s = $c.Style
$c.Value = "'113.25"
$c.Style = $s
$ws.Range("E50").Value = '  +1.82%  '
$c = $ws.Range("D51")
$s = $c.Style
$c.Value = "'0.2967"
$c.Style = $s
$ws.Range("E51").Value = '  +10.17%  '
